# aggiornamento fino a 1/09/2021
# Append 9 new daily rows (r358:r366) to Sheet1, covering 2021-08-24 .. 2021-09-01
# (Excel serials 44432..44440). New rows get the same formatting as the last
# existing row (357) by copying it down first, then overwriting the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting (style/number format/borders) of the last data row into
# the new block of rows.
$ws.Range("A357:D357").Copy($ws.Range("A358:D366"))

$ws.Cells.Item(358,1).Value = 44432
$ws.Cells.Item(358,2).Value = 1
$ws.Cells.Item(358,3).Value = 2
$ws.Cells.Item(358,4).Value = 87.56567425569177

$ws.Cells.Item(359,1).Value = 44433
$ws.Cells.Item(359,2).Value = 0
$ws.Cells.Item(359,3).Value = 2
$ws.Cells.Item(359,4).Value = 87.56567425569177

$ws.Cells.Item(360,1).Value = 44434
$ws.Cells.Item(360,2).Value = 0
$ws.Cells.Item(360,3).Value = 2
$ws.Cells.Item(360,4).Value = 87.56567425569177

$ws.Cells.Item(361,1).Value = 44435
$ws.Cells.Item(361,2).Value = 0
$ws.Cells.Item(361,3).Value = 2
$ws.Cells.Item(361,4).Value = 87.56567425569177

$ws.Cells.Item(362,1).Value = 44436
$ws.Cells.Item(362,2).Value = 0
$ws.Cells.Item(362,3).Value = 2
$ws.Cells.Item(362,4).Value = 87.56567425569177

$ws.Cells.Item(363,1).Value = 44437
$ws.Cells.Item(363,2).Value = 0
$ws.Cells.Item(363,3).Value = 1
$ws.Cells.Item(363,4).Value = 43.78283712784589

$ws.Cells.Item(364,1).Value = 44438
$ws.Cells.Item(364,2).Value = 2
$ws.Cells.Item(364,3).Value = 3
$ws.Cells.Item(364,4).Value = 131.3485113835376

$ws.Cells.Item(365,1).Value = 44439
$ws.Cells.Item(365,2).Value = 1
$ws.Cells.Item(365,3).Value = 3
$ws.Cells.Item(365,4).Value = 131.3485113835376

$ws.Cells.Item(366,1).Value = 44440
$ws.Cells.Item(366,2).Value = 0
$ws.Cells.Item(366,3).Value = 3
$ws.Cells.Item(366,4).Value = 131.3485113835376
